$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(19, 8).Value = 2287.8572
$ws.Cells.Item(19, 9).Value = 1900
$ws.Cells.Item(19, 10).Value = 2443
$ws.Cells.Item(19, 11).Value = 1900
$ws.Cells.Item(19, 12).Value = 2443
$ws.Cells.Item(19, 13).Value = -1725
$ws.Cells.Item(19, 14).Value = -2793

$ws.Cells.Item(75, 8).Value = 23438
$ws.Cells.Item(75, 9).Value = 20000
$ws.Cells.Item(75, 10).Value = 25157
$ws.Cells.Item(75, 11).Value = 20000
$ws.Cells.Item(75, 12).Value = 25157
$ws.Cells.Item(75, 13).Value = -19064
$ws.Cells.Item(75, 14).Value = -27029

$ws.Cells.Item(78, 8).Value = 23438
$ws.Cells.Item(78, 9).Value = 20000
$ws.Cells.Item(78, 10).Value = 25157
$ws.Cells.Item(78, 11).Value = 60000
$ws.Cells.Item(78, 12).Value = 75471
$ws.Cells.Item(78, 13).Value = -55320
$ws.Cells.Item(78, 14).Value = -84831

$ws.Cells.Item(93, 8).Value = 0
$ws.Cells.Item(93, 9).Value = 0
$ws.Cells.Item(93, 10).Value = 0
$ws.Cells.Item(93, 11).Value = 0
$ws.Cells.Item(93, 12).Value = 0
$ws.Cells.Item(93, 14).ClearContents()

$ws.Cells.Item(137, 8).Value = 16394585
$ws.Cells.Item(137, 9).Value = 1018.1591
$ws.Cells.Item(137, 10).Value = 58824990
$ws.Cells.Item(137, 11).Value = 3054.4773
$ws.Cells.Item(137, 12).Value = 176474970
$ws.Cells.Item(137, 13).Value = -504.4773
$ws.Cells.Item(137, 14).Value = -176480070

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 920241
$ws.Cells.Item(2, 9).Value = 863.1579
$ws.Cells.Item(2, 10).Value = 2263947.2
$ws.Cells.Item(2, 11).Value = 863.1579
$ws.Cells.Item(2, 12).Value = 2263947.2
$ws.Cells.Item(2, 13).Value = -750.1579
$ws.Cells.Item(2, 14).Value = -2264173.2

$ws.Cells.Item(43, 8).Value = 11555
$ws.Cells.Item(43, 9).Value = 0
$ws.Cells.Item(43, 10).Value = 11555
$ws.Cells.Item(43, 11).Value = 0
$ws.Cells.Item(43, 12).Value = 11555
$ws.Cells.Item(43, 14).Value = -12181

$ws.Cells.Item(61, 8).Value = 1544034
$ws.Cells.Item(61, 9).Value = 1764508.1
$ws.Cells.Item(61, 10).Value = 715.55554
$ws.Cells.Item(61, 11).Value = 1764508.1
$ws.Cells.Item(61, 12).Value = 715.55554
$ws.Cells.Item(61, 13).Value = -1764296.1
$ws.Cells.Item(61, 14).Value = -1139.55554

$ws.Cells.Item(94, 8).Value = 15000
$ws.Cells.Item(94, 9).Value = 0
$ws.Cells.Item(94, 10).Value = 15000
$ws.Cells.Item(94, 11).Value = 0
$ws.Cells.Item(94, 12).Value = 15000
$ws.Cells.Item(94, 13).ClearContents()
$ws.Cells.Item(94, 14).Value = -16802

$ws.Cells.Item(116, 8).Value = 920241
$ws.Cells.Item(116, 9).Value = 863.1579
$ws.Cells.Item(116, 10).Value = 2263947.2
$ws.Cells.Item(116, 11).Value = 863.1579
$ws.Cells.Item(116, 12).Value = 2263947.2
$ws.Cells.Item(116, 13).Value = 1430.8421
$ws.Cells.Item(116, 14).Value = -2268535.2

$ws.Cells.Item(132, 8).Value = 4923865.5
$ws.Cells.Item(132, 9).Value = 5612796
$ws.Cells.Item(132, 10).Value = 101350.9
$ws.Cells.Item(132, 11).Value = 16838388
$ws.Cells.Item(132, 12).Value = 304052.7
$ws.Cells.Item(132, 13).Value = -16835858
$ws.Cells.Item(132, 14).Value = -309112.7

$ws.Cells.Item(136, 8).Value = 1544034
$ws.Cells.Item(136, 9).Value = 1764508.1
$ws.Cells.Item(136, 10).Value = 715.55554
$ws.Cells.Item(136, 11).Value = 5293524.300000001
$ws.Cells.Item(136, 12).Value = 2146.66662
$ws.Cells.Item(136, 13).Value = -5290974.300000001
$ws.Cells.Item(136, 14).Value = -7246.66662

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 920241
$ws.Cells.Item(3, 9).Value = 863.1579
$ws.Cells.Item(3, 10).Value = 2263947.2
$ws.Cells.Item(3, 11).Value = 863.1579
$ws.Cells.Item(3, 12).Value = 2263947.2
$ws.Cells.Item(3, 13).Value = -749.1579
$ws.Cells.Item(3, 14).Value = -2264175.2

$ws.Cells.Item(76, 8).Value = 30000
$ws.Cells.Item(76, 9).Value = 0
$ws.Cells.Item(76, 10).Value = 30000
$ws.Cells.Item(76, 11).Value = 0
$ws.Cells.Item(76, 12).Value = 30000
$ws.Cells.Item(76, 14).Value = -30630

$ws.Cells.Item(79, 8).Value = 30000
$ws.Cells.Item(79, 9).Value = 0
$ws.Cells.Item(79, 10).Value = 30000
$ws.Cells.Item(79, 11).Value = 0
$ws.Cells.Item(79, 12).Value = 30000
$ws.Cells.Item(79, 14).Value = -32184

$ws.Cells.Item(94, 8).Value = 726.625
$ws.Cells.Item(94, 9).Value = 636.95
$ws.Cells.Item(94, 10).Value = 1175
$ws.Cells.Item(94, 11).Value = 636.95
$ws.Cells.Item(94, 12).Value = 1175
$ws.Cells.Item(94, 13).Value = -185.95
$ws.Cells.Item(94, 14).Value = -2077

$ws.Cells.Item(134, 8).Value = 4279539
$ws.Cells.Item(134, 9).Value = 4788865
$ws.Cells.Item(134, 10).Value = 1200
$ws.Cells.Item(134, 11).Value = 14366595
$ws.Cells.Item(134, 12).Value = 3600
$ws.Cells.Item(134, 13).Value = -14364060
$ws.Cells.Item(134, 14).Value = -8670

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(5, 8).Value = 400.3684
$ws.Cells.Item(5, 9).Value = 140.85715
$ws.Cells.Item(5, 10).Value = 551.75
$ws.Cells.Item(5, 11).Value = 140.85715
$ws.Cells.Item(5, 12).Value = 551.75
$ws.Cells.Item(5, 13).Value = -28.85714999999999
$ws.Cells.Item(5, 14).Value = -775.75

$ws.Cells.Item(50, 8).Value = 10127.223
$ws.Cells.Item(50, 9).Value = 7781
$ws.Cells.Item(50, 10).Value = 11300.333
$ws.Cells.Item(50, 11).Value = 7781
$ws.Cells.Item(50, 12).Value = 11300.333
$ws.Cells.Item(50, 13).Value = -7156
$ws.Cells.Item(50, 14).Value = -12550.333

$ws.Cells.Item(53, 8).Value = 33995
$ws.Cells.Item(53, 9).Value = 0
$ws.Cells.Item(53, 10).Value = 33995
$ws.Cells.Item(53, 11).Value = 0
$ws.Cells.Item(53, 12).Value = 33995
$ws.Cells.Item(53, 14).Value = -35209

$ws.Cells.Item(58, 8).Value = 1577.716
$ws.Cells.Item(58, 9).Value = 748.4643
$ws.Cells.Item(58, 10).Value = 3028.9062
$ws.Cells.Item(58, 11).Value = 748.4643
$ws.Cells.Item(58, 12).Value = 3028.9062
$ws.Cells.Item(58, 13).Value = -545.4643
$ws.Cells.Item(58, 14).Value = -3434.9062

$ws.Cells.Item(59, 8).Value = 15799.637
$ws.Cells.Item(59, 9).Value = 10000
$ws.Cells.Item(59, 10).Value = 16379.6
$ws.Cells.Item(59, 11).Value = 10000
$ws.Cells.Item(59, 12).Value = 16379.6
$ws.Cells.Item(59, 13).Value = -8855
$ws.Cells.Item(59, 14).Value = -18669.6

$ws.Cells.Item(60, 8).Value = 8753.272000000001
$ws.Cells.Item(60, 9).Value = 4993.3335
$ws.Cells.Item(60, 10).Value = 10163.25
$ws.Cells.Item(60, 11).Value = 4993.3335
$ws.Cells.Item(60, 12).Value = 10163.25
$ws.Cells.Item(60, 13).Value = -4482.3335
$ws.Cells.Item(60, 14).Value = -11185.25

$ws.Cells.Item(68, 8).Value = 17114.572
$ws.Cells.Item(68, 9).Value = 10000
$ws.Cells.Item(68, 10).Value = 19960.4
$ws.Cells.Item(68, 11).Value = 10000
$ws.Cells.Item(68, 12).Value = 19960.4
$ws.Cells.Item(68, 13).Value = -9251
$ws.Cells.Item(68, 14).Value = -21458.4

$ws.Cells.Item(71, 8).Value = 17114.572
$ws.Cells.Item(71, 9).Value = 10000
$ws.Cells.Item(71, 10).Value = 19960.4
$ws.Cells.Item(71, 11).Value = 30000
$ws.Cells.Item(71, 12).Value = 59881.2
$ws.Cells.Item(71, 13).Value = -26256
$ws.Cells.Item(71, 14).Value = -67369.20000000001

$ws.Cells.Item(74, 8).Value = 13135.223
$ws.Cells.Item(74, 9).Value = 2000
$ws.Cells.Item(74, 10).Value = 14527.125
$ws.Cells.Item(74, 11).Value = 2000
$ws.Cells.Item(74, 12).Value = 14527.125
$ws.Cells.Item(74, 13).Value = -1126
$ws.Cells.Item(74, 14).Value = -16275.125

$ws.Cells.Item(77, 8).Value = 13135.223
$ws.Cells.Item(77, 9).Value = 2000
$ws.Cells.Item(77, 10).Value = 14527.125
$ws.Cells.Item(77, 11).Value = 6000
$ws.Cells.Item(77, 12).Value = 43581.375
$ws.Cells.Item(77, 13).Value = -1632
$ws.Cells.Item(77, 14).Value = -52317.375

$ws.Cells.Item(132, 8).Value = 1376.1147
$ws.Cells.Item(132, 9).Value = 1265.4386
$ws.Cells.Item(132, 10).Value = 2953.25
$ws.Cells.Item(132, 11).Value = 3796.3158
$ws.Cells.Item(132, 12).Value = 8859.75
$ws.Cells.Item(132, 13).Value = -1266.3158
$ws.Cells.Item(132, 14).Value = -13919.75

$ws.Cells.Item(136, 8).Value = 1577.716
$ws.Cells.Item(136, 9).Value = 748.4643
$ws.Cells.Item(136, 10).Value = 3028.9062
$ws.Cells.Item(136, 11).Value = 2245.3929
$ws.Cells.Item(136, 12).Value = 9086.7186
$ws.Cells.Item(136, 13).Value = 304.6071000000002
$ws.Cells.Item(136, 14).Value = -14186.7186

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(105, 8).Value = 2002249.5
$ws.Cells.Item(105, 9).Value = 0
$ws.Cells.Item(105, 10).Value = 2002249.5
$ws.Cells.Item(105, 11).Value = 0
$ws.Cells.Item(105, 12).Value = 6006748.5
$ws.Cells.Item(105, 14).Value = -6011990.5

$ws.Cells.Item(110, 8).Value = 2857.5
$ws.Cells.Item(110, 9).Value = 2450
$ws.Cells.Item(110, 10).Value = 3265
$ws.Cells.Item(110, 11).Value = 7350
$ws.Cells.Item(110, 12).Value = 9795
$ws.Cells.Item(110, 13).Value = -3260
$ws.Cells.Item(110, 14).Value = -17975

$ws.Cells.Item(113, 8).Value = 452.35483
$ws.Cells.Item(113, 9).Value = 440.3
$ws.Cells.Item(113, 10).Value = 474.27274
$ws.Cells.Item(113, 11).Value = 1320.9
$ws.Cells.Item(113, 12).Value = 1422.81822
$ws.Cells.Item(113, 13).Value = 849.0999999999999
$ws.Cells.Item(113, 14).Value = -5762.81822

$ws.Cells.Item(131, 8).Value = 3512.878
$ws.Cells.Item(131, 9).Value = 6497
$ws.Cells.Item(131, 10).Value = 2550.258
$ws.Cells.Item(131, 11).Value = 19491
$ws.Cells.Item(131, 12).Value = 7650.773999999999
$ws.Cells.Item(131, 13).Value = -14451
$ws.Cells.Item(131, 14).Value = -17730.774

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(10, 8).Value = 70004.5
$ws.Cells.Item(10, 9).Value = 0
$ws.Cells.Item(10, 10).Value = 70004.5
$ws.Cells.Item(10, 11).Value = 0
$ws.Cells.Item(10, 12).Value = 70004.5
$ws.Cells.Item(10, 14).Value = -70342.5

$ws.Cells.Item(80, 8).Value = 9650.5
$ws.Cells.Item(80, 9).Value = 0
$ws.Cells.Item(80, 10).Value = 9650.5
$ws.Cells.Item(80, 11).Value = 0
$ws.Cells.Item(80, 12).Value = 9650.5
$ws.Cells.Item(80, 14).Value = -11646.5

$ws.Cells.Item(83, 8).Value = 9650.5
$ws.Cells.Item(83, 9).Value = 0
$ws.Cells.Item(83, 10).Value = 9650.5
$ws.Cells.Item(83, 11).Value = 0
$ws.Cells.Item(83, 12).Value = 28951.5
$ws.Cells.Item(83, 14).Value = -38935.5

$ws.Cells.Item(132, 8).Value = 7692926
$ws.Cells.Item(132, 9).Value = 9569093
$ws.Cells.Item(132, 10).Value = 641
$ws.Cells.Item(132, 11).Value = 28707279
$ws.Cells.Item(132, 12).Value = 1923
$ws.Cells.Item(132, 13).Value = -28704749
$ws.Cells.Item(132, 14).Value = -6983
